# Update "想去人数" (wanted-to-go count) figures for two 南宁 events
# that appear on both the "展览" and "全部类型" sheets.
#   展览:    F3 1125 -> 1127 ; F4 2562 -> 2563
#   全部类型: F5 1125 -> 1127 ; F6 2562 -> 2563

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1127
$wsExhibit.Range("F4").Value = 2563

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1127
$wsAll.Range("F6").Value = 2563
